$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3696246258455235
$ws.Range("J2").Value = 0.2702026041520829
$ws.Range("K2").Value = -0.2461390038637574
$ws.Range("L2").Value = 2.549338570627156

$ws.Range("I20").Value = 0.1617039438647276
$ws.Range("J20").Value = 0.3695633939073046
$ws.Range("K20").Value = 0.01276248446693575
$ws.Range("L20").Value = 2.004878396784024
